$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove row 10 ("data/simulated_u.csv and data/simulated_v.csv" / "dataset for
# Example 3"). This shifts every subsequent row up by one, which is exactly
# what the target workbook shows (old row 11 becomes new row 10, etc.).
$ws.Rows.Item(10).Delete()

# After the shift, update the description text (column E) for a handful of
# rows to match the new wording used in the target workbook.
$ws.Range("E13").Value = "reproduce results in subsection 3.3: draw two types of local Kendall's tau surfaces"
$ws.Range("E16").Value = "reproduce simulation results in Section 5"
$ws.Range("E18").Value = "generate all plots in Figure S.16 of the Supplementary Materials "
$ws.Range("E19").Value = "data for Figure S.16"
$ws.Range("E21").Value = "reproduce simulation results in Section 4: comparative analysis of global and local kendall's tau for copula models"
$ws.Range("E22").Value = "reproduce results in Figures S.7, S.9, S.11, and S.13"
$ws.Range("E23").Value = "reproduce results in Figures S.8, S.10, S.12, and S.14"

# Match the saved selection state of the target workbook.
$ws.Range("C26").Select()
